$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Depot API changed with cleaner structure: the SBU/category codes for the
# Bandarban depot field-colleague rows (8-13) are now distinct per field
# colleague instead of reusing "A"/"B" for most rows.
$ws.Range("F8").Value = "C"
$ws.Range("F9").Value = "D"
$ws.Range("F10").Value = "E"
$ws.Range("F11").Value = "N"
$ws.Range("F12").Value = "B"
$ws.Range("F13").Value = "M"

# Leave the selection where the edit ended, on the last touched cell.
$ws.Range("F13").Select()
